$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency number format used by existing "Valor"/"Preço unitário" columns
$currencyFmt = '_-"R$"* #,##0.00######_-;"-R$"* #,##0.00######_-;_ \-;_-@_-'

# ---------------------------------------------------------------------
# 1) Insert two new rows before row 6 for the PETR4 Dividendo / JCP
#    entries (pushes the old rows 6-9 down to 8-11).
# ---------------------------------------------------------------------
$ws.Range("A6:A7").EntireRow.Insert()

# --- Row 6: Credito / 15/12/2022 / Dividendo / PETR4.../ Corretora A / 2 / 2.2 / 4.39
$ws.Range("A6").Value = "Credito"
$ws.Range("B6").Value = "15/12/2022"
$ws.Range("C6").Value = "Dividendo"
$ws.Range("D6").Value = "PETR4 - PETROLEO BRASILEIRO S/A PETROBRAS"
$ws.Range("E6").Value = "Corretora A"

# --- Row 7: Credito / 15/12/2022 / Juros Sobre Capital Próprio / PETR4.../ Corretora A / 2 / 1.06 / 1.8
$ws.Range("A7").Value = "Credito"
$ws.Range("B7").Value = "15/12/2022"
$ws.Range("C7").Value = "Juros Sobre Capital Próprio"
$ws.Range("D7").Value = "PETR4 - PETROLEO BRASILEIRO S/A PETROBRAS"
$ws.Range("E7").Value = "Corretora A"

# Column F ("Quantidade") stores the literal text "2" on these two rows
# (not a number) -- reset the cells to the plain/default style first
# (copy A6's style, which is the bare default), flip to text format,
# write the string, then restore the plain style.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("F6:F7").PasteSpecial(-4122) | Out-Null
$ws.Range("F6:F7").NumberFormat = "@"
$ws.Range("F6").Value = "2"
$ws.Range("F7").Value = "2"
$ws.Range("A6").Copy() | Out-Null
$ws.Range("F6:F7").PasteSpecial(-4122) | Out-Null

$ws.Range("G6:H7").NumberFormat = $currencyFmt
$ws.Range("G6").Value = 2.2
$ws.Range("H6").Value = 4.39
$ws.Range("G7").Value = 1.06
$ws.Range("H7").Value = 1.8

$ws.Rows.Item(6).RowHeight = 13.4
$ws.Rows.Item(7).RowHeight = 13.4

# ---------------------------------------------------------------------
# 2) Append two new rows (12 and 13) with the ITSA4 JCP / Dividendo
#    entries.
# ---------------------------------------------------------------------
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B12:B13").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Value = "Credito"
$ws.Range("B12").Value = "26/08/2022"
$ws.Range("C12").Value = "Juros Sobre Capital Próprio"
$ws.Range("D12").Value = "ITSA4 - ITAUSA S/A"
$ws.Range("E12").Value = "Corretora A"

$ws.Range("A13").Value = "Credito"
$ws.Range("B13").Value = 44799
$ws.Range("C13").Value = "Dividendo"
$ws.Range("D13").Value = "ITSA4 - ITAUSA S.A.                                       "
$ws.Range("E13").Value = "Corretora A"

$ws.Range("A11").Copy() | Out-Null
$ws.Range("F12:F13").PasteSpecial(-4122) | Out-Null
$ws.Range("F12:F13").NumberFormat = "@"
$ws.Range("F12").Value = "16"
$ws.Range("F13").Value = "16"
$ws.Range("A11").Copy() | Out-Null
$ws.Range("F12:F13").PasteSpecial(-4122) | Out-Null

$ws.Range("G12:H13").NumberFormat = $currencyFmt
$ws.Range("G12").Value = 0.04
$ws.Range("H12").Value = 0.51
$ws.Range("G13").Value = 0.02
$ws.Range("H13").Value = 0.32

$ws.Rows.Item(12).RowHeight = 13.4

# ---------------------------------------------------------------------
# 3) Misc sheet-level bookkeeping to match the edited workbook state.
# ---------------------------------------------------------------------
$ws.Range("B17").Select() | Out-Null
